# Add two new date columns (T = "T18: 6/4/2020", U = "T19: 7/4/2020") to the
# COVID-19 history table, carrying forward / updating the per-department case
# counts, and recomputing the "Sum" row totals for the two new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting of column S (the previous last data column) onto
#        the two new columns T and U, row by row, so each cell keeps the same
#        style as its neighbour in column S. ------------------------------
for ($r = 1; $r -le 20; $r++) {
    $ws.Range("S$r").Copy()
    $ws.Range("T$r").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("S$r").Copy()
    $ws.Range("U$r").PasteSpecial(-4122)   # xlPasteFormats
}

# --- 2. Header row: new date labels -----------------------------------------
$ws.Range("T1").Value = "T18: 6/4/2020"
$ws.Range("U1").Value = "T19: 7/4/2020"

# --- 3. Per-department data, rows 2-19 (Atlantida ... Yoro) -----------------
$tvals = @(11, 2, 20, 2, 1, 195, 0, 53, 0, 0, 0, 1, 4, 0, 0, 8, 0, 8)
$uvals = @(16, 2, 20, 2, 1, 195, 0, 54, 0, 0, 0, 1, 4, 0, 0, 8, 0, 9)

for ($i = 0; $i -lt 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $tvals[$i]
    $ws.Cells.Item($row, 21).Value = $uvals[$i]
}

# --- 4. Sum row (row 20) -----------------------------------------------------
$ws.Range("T20").Formula = "=SUM(T2:T19)"
$ws.Range("U20").Formula = "=SUM(U2:U19)"

# --- 5. Column widths for the two new columns --------------------------------
# (XLSX stores width in "characters"; ColumnWidth set via COM is offset from
# the stored value by the sheet's max-digit-width quantum, so these inputs
# are chosen to land as close as possible to the target stored widths of
# 14.86 and 15.)
$ws.Range("T1").EntireColumn.ColumnWidth = 14.0
$ws.Range("U1").EntireColumn.ColumnWidth = 14.1667

# --- 6. Selection / active cell ----------------------------------------------
$ws.Range("U20").Select()

Write-Host "Added columns T (T18: 6/4/2020) and U (T19: 7/4/2020)"
